$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.404.62"
$ws.Range("E2").Value = "  -1.89%  "

$ws.Range("D3").Value = "'1.905.48"
$ws.Range("E3").Value = "  -2.55%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'238.08"
$ws.Range("E5").Value = "  -2.41%  "

$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").Value = "'0.4721"
$ws.Range("E7").Value = "  -2.50%  "

$ws.Range("D8").Value = "'0.2828"
$ws.Range("E8").Value = "  -3.85%  "

$ws.Range("D9").Value = "'0.06629"
$ws.Range("E9").Value = "  -6.48%  "

$ws.Range("E10").Value = "  -6.53%  "

$ws.Range("D11").Value = "'99.92"
$ws.Range("E11").Value = "  -6.81%  "

$ws.Range("D12").Value = "'0.07704"
$ws.Range("E12").Value = "  -1.10%  "

$ws.Range("D13").Value = "'1.913.65"
$ws.Range("E13").Value = "  -2.15%  "

$ws.Range("D14").Value = "'5.169"
$ws.Range("E14").Value = "  -5.00%  "

$ws.Range("D15").Value = "'0.6656"
$ws.Range("E15").Value = "  -5.06%  "

$ws.Range("D16").Value = "'30.373.09"
$ws.Range("E16").Value = "  -2.02%  "

$ws.Range("D17").Value = "'253.80"
$ws.Range("E17").Value = "  -9.35%  "

$ws.Range("E18").Value = "  +0.08%  "

$ws.Range("D19").Value = "'0.000007432"
$ws.Range("E19").Value = "  -5.07%  "

$ws.Range("D20").Value = "'12.60"
$ws.Range("E20").Value = "  -5.12%  "

$ws.Range("D21").Value = "'5.342"
$ws.Range("E21").Value = "  -3.71%  "

$ws.Range("D22").Value = "'1.003"
$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "'6.273"
$ws.Range("E23").Value = "  -3.58%  "

$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "'9.318"
$ws.Range("E24").Value = "  -5.01%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'165.42"
$ws.Range("E25").Value = "  -2.10%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'18.82"
$ws.Range("E26").Value = "  -4.59%  "

$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'2.032"
$ws.Range("E27").Value = "  -6.52%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'1.383"
$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("D29").Value = "'0.1006"
$ws.Range("E29").Value = "  -3.93%  "

$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'4.633"
$ws.Range("E30").Value = "  +0.55%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.506"
$ws.Range("E31").Value = "  -4.11%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'4.216"
$ws.Range("E32").Value = "  -5.19%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.04705"
$ws.Range("E33").Value = "  -3.96%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7226"
$ws.Range("E34").Value = "  -3.24%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.100"
$ws.Range("E35").Value = "  -5.69%  "

$ws.Range("B36").Value = "Frax"
$ws.Range("C36").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D36").Value = "'1.0000"
$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.716"
$ws.Range("E37").Value = "  -0.70%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01894"
$ws.Range("E38").Value = "  -5.20%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.595"
$ws.Range("E39").Value = "  -3.48%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'6.190"
$ws.Range("E40").Value = "  -5.32%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'72.43"
$ws.Range("E41").Value = "  -6.80%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'1.969"
$ws.Range("E42").Value = "  -7.16%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'106.10"
$ws.Range("E43").Value = "  -2.86%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8551"
$ws.Range("E44").Value = "  -4.90%  "

$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'1.001"
$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").Value = "'0.4211"
$ws.Range("E46").Value = "  -5.29%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "'1.024.93"
$ws.Range("E47").Value = "  +3.32%  "

$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.360"
$ws.Range("E48").Value = "  -7.63%  "

$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.1188"
$ws.Range("E49").Value = "  -4.61%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'34.40"
$ws.Range("E50").Value = "  -4.07%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'8.708"
$ws.Range("E51").Value = "  -6.36%  "
